# Merge the two trailing runs of the 3rd paragraph in "Retângulo 40"
# (slide 10, shape 11) back into a single run, i.e.
#   " nesse tal lugar, que uma vida não chega para lá " + "chegar."
# becomes
#   " nesse tal lugar, que uma vida não chega para lá chegar."
# (both runs already share identical character formatting).

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(10)
$shp = $s.Shapes.Item(11)

# The shape auto-fits its height/width to the text ("spAutoFit"); this
# COM shim recomputes that box as soon as any text is written, so stash
# the current size (read *before* touching the text) and restore it
# afterwards. The Height/Width getters round-trip the underlying EMU
# value through a 32-bit float, which truncates by up to ~1 EMU on the
# way back in; nudging by half an EMU (in points) before writing
# compensates for that truncation so the restored size lands back on
# the exact original EMU value.
$origHeight = $shp.Height
$origWidth  = $shp.Width
$epsilon    = 0.5 / 12700.0

$tr    = $shp.TextFrame.TextRange
$para3 = $tr.Paragraphs(3)

$fullText = $para3.Text
$startIdx = $fullText.IndexOf(" nesse tal lugar")
$rng = $para3.Characters($startIdx + 1, $fullText.Length - $startIdx)
$rng.Text = " nesse tal lugar, que uma vida não chega para lá chegar."

$shp.Height = $origHeight + $epsilon
$shp.Width  = $origWidth + $epsilon
